# Lattice multiplication exercise sheet update.
# The worksheet table (5 rows x 3 columns = 15 cells) gets every problem
# replaced with a new "N x M" pair together with its matching lattice grid
# labels. Each cell's text is: "A x B", the tens/ones digits of B, a dashed
# divider, and two rows showing the tens/ones digits of A down the left
# edge of the lattice grid.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$vb = [char]11   # vertical-tab == <w:br/> line break inside a Word run

# Row, Col, "A x B" text, "  tensB    onesB", tensA left-label line, onesA left-label line
$cells = @(
    @(1, 1, "87 x 74", "  7    4", "8|    |", "7|    |"),
    @(1, 2, "30 x 28", "  2    8", "3|    |", "0|    |"),
    @(1, 3, "84 x 12", "  1    2", "8|    |", "4|    |"),

    @(2, 1, "63 x 13", "  1    3", "6|    |", "3|    |"),
    @(2, 2, "84 x 63", "  6    3", "8|    |", "4|    |"),
    @(2, 3, "27 x 86", "  8    6", "2|    |", "7|    |"),

    @(3, 1, "57 x 29", "  2    9", "5|    |", "7|    |"),
    @(3, 2, "75 x 76", "  7    6", "7|    |", "5|    |"),
    @(3, 3, "71 x 87", "  8    7", "7|    |", "1|    |"),

    @(4, 1, "11 x 40", "  4    0", "1|    |", "1|    |"),
    @(4, 2, "78 x 94", "  9    4", "7|    |", "8|    |"),
    @(4, 3, "54 x 60", "  6    0", "5|    |", "4|    |"),

    @(5, 1, "56 x 93", "  9    3", "5|    |", "6|    |"),
    @(5, 2, "72 x 43", "  4    3", "7|    |", "2|    |"),
    @(5, 3, "75 x 76", "  7    6", "7|    |", "5|    |")
)

foreach ($item in $cells) {
    $row = $item[0]
    $col = $item[1]
    $line1 = $item[2]
    $line2 = $item[3]
    $line4 = $item[4]
    $line5 = $item[5]

    $newText = $line1 + $vb + $line2 + $vb + "  ----" + $vb + $line4 + $vb + $line5
    $table.Cell($row, $col).Range.Text = $newText
}
